$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[João Paulo-Sistemas digitais-2A, João Paulo-Sistemas digitais-2A, -]"
$ws.Range("F2").Value = "-"

$ws.Range("B3").Value = "[André Guimarães-CAD-2A, André Guimarães-CAD-2A]"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "Andre B.-Eletrônica analóg. e de potência"
$ws.Range("E3").Value = "[Pedro Bispo-Acionamentos Elétricos-2A, -, -, Sandro-Programação de Computadores-2A]"
$ws.Range("F3").Value = "-"

$ws.Range("B4").Value = "Sandro-Circuitos elétricos 2"
$ws.Range("D4").Value = "Andre B.-Eletrônica analóg. e de potência"
$ws.Range("E4").Value = "[Pedro Bispo-Acionamentos Elétricos-2A, -, -, Sandro-Programação de Computadores-2A]"
$ws.Range("F4").Value = "-"

$ws.Range("B6").Value = "Sandro-Circuitos elétricos 2"
$ws.Range("E6").Value = "[Pedro Bispo-Acionamentos Elétricos-2A, -, -, Sandro-Programação de Computadores-2A]"
$ws.Range("F6").Value = "-"

$ws.Range("B7").Value = "Nilton Maia-M.T.R.M."
$ws.Range("E7").Value = "[Pedro Bispo-Acionamentos Elétricos-2A, -, -, Sandro-Programação de Computadores-2A]"
$ws.Range("F7").Value = "-"

$ws.Range("B8").Value = "Nilton Maia-M.T.R.M."
$ws.Range("F8").Value = "-"

$ws.Range("D15").Value = "-"
